$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (changed) date, stored as a date serial.
# Every populated row (2-29) moves from serial 45562 (2024-09-27) to
# 45563 (2024-09-28) -- a routine "last updated" bump.
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45562) {
        $cell.Value = 45563
    }
}
